# Auto-generated Excel COM-interop script to apply scheduled market-data update
# to the Cactuar_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 669.44446
$ws.Range("I9").Value = 654.1667
$ws.Range("K9").Value = 654.1667
$ws.Range("M9").Value = -485.1667
$ws.Range("H34").Value = 19497
$ws.Range("I34").Value = 19497
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 19497
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -19294
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 19497
$ws.Range("I36").Value = 19497
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 19497
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -18782
$ws.Range("N36").ClearContents()
$ws.Range("H38").Value = 6405.125
$ws.Range("I38").Value = 3149.3333
$ws.Range("J38").Value = 8358.6
$ws.Range("K38").Value = 9447.999899999999
$ws.Range("L38").Value = 25075.8
$ws.Range("M38").Value = -9075.999899999999
$ws.Range("N38").Value = -25819.8
$ws.Range("H40").Value = 18532604
$ws.Range("I40").Value = 15047
$ws.Range("J40").Value = 45467236
$ws.Range("K40").Value = 15047
$ws.Range("L40").Value = 45467236
$ws.Range("M40").Value = -14872
$ws.Range("N40").Value = -45467586
$ws.Range("H51").Value = 5249.7144
$ws.Range("J51").Value = 5130.75
$ws.Range("L51").Value = 5130.75
$ws.Range("N51").Value = -6098.75
$ws.Range("H100").Value = 1694.5
$ws.Range("I100").Value = 1694.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1694.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1153.5
$ws.Range("N100").ClearContents()
$ws.Range("H121").Value = 4371.2
$ws.Range("J121").Value = 4371.2
$ws.Range("L121").Value = 13113.6
$ws.Range("N121").Value = -16607.6
$ws.Range("H129").Value = 1434.0667
$ws.Range("I129").Value = 865.4
$ws.Range("K129").Value = 2596.2
$ws.Range("M129").Value = 2403.8
$ws.Range("H138").Value = 3919.9656
$ws.Range("I138").Value = 3082.3
$ws.Range("J138").Value = 4360.8423
$ws.Range("K138").Value = 9246.900000000001
$ws.Range("L138").Value = 13082.5269
$ws.Range("M138").Value = -4106.900000000001
$ws.Range("N138").Value = -23362.5269

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13661.243
$ws.Range("I32").Value = 14006.018
$ws.Range("K32").Value = 14006.018
$ws.Range("M32").Value = -13719.018
$ws.Range("H45").Value = 4539.4165
$ws.Range("I45").Value = 4135.8887
$ws.Range("K45").Value = 4135.8887
$ws.Range("M45").Value = -3758.8887
$ws.Range("H46").Value = 12474.286
$ws.Range("J46").Value = 12664.615
$ws.Range("L46").Value = 12664.615
$ws.Range("N46").Value = -13302.615
$ws.Range("H63").Value = 1934.5
$ws.Range("I63").Value = 1860.4
$ws.Range("K63").Value = 1860.4
$ws.Range("M63").Value = -1174.4
$ws.Range("H66").Value = 1934.5
$ws.Range("I66").Value = 1860.4
$ws.Range("K66").Value = 9302
$ws.Range("M66").Value = -5870
$ws.Range("H74").Value = 5320014
$ws.Range("I74").Value = 7143396.5
$ws.Range("J74").Value = 1815.3334
$ws.Range("K74").Value = 7143396.5
$ws.Range("L74").Value = 1815.3334
$ws.Range("M74").Value = -7142522.5
$ws.Range("N74").Value = -3563.3334
$ws.Range("H77").Value = 5320014
$ws.Range("I77").Value = 7143396.5
$ws.Range("J77").Value = 1815.3334
$ws.Range("K77").Value = 35716982.5
$ws.Range("L77").Value = 9076.666999999999
$ws.Range("M77").Value = -35712614.5
$ws.Range("N77").Value = -17812.667
$ws.Range("H122").Value = 4384
$ws.Range("I122").Value = 3462.4348
$ws.Range("J122").Value = 7916.6665
$ws.Range("K122").Value = 10387.3044
$ws.Range("L122").Value = 23749.9995
$ws.Range("M122").Value = -7937.304400000001
$ws.Range("N122").Value = -28649.9995
$ws.Range("H132").Value = 14531.108
$ws.Range("I132").Value = 16917.527
$ws.Range("K132").Value = 50752.58099999999
$ws.Range("M132").Value = -48222.58099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 623360.3
$ws.Range("I94").Value = 685496.5
$ws.Range("K94").Value = 685496.5
$ws.Range("M94").Value = -685045.5
$ws.Range("H134").Value = 4817.143
$ws.Range("I134").Value = 1683.125
$ws.Range("K134").Value = 5049.375
$ws.Range("M134").Value = -2514.375
$ws.Range("H138").Value = 97593
$ws.Range("J138").Value = 97593
$ws.Range("L138").Value = 97593
$ws.Range("N138").Value = -107873

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5804.114
$ws.Range("I31").Value = 1437.5834
$ws.Range("K31").Value = 1437.5834
$ws.Range("M31").Value = -1142.5834
$ws.Range("H34").Value = 5804.114
$ws.Range("I34").Value = 1437.5834
$ws.Range("K34").Value = 1437.5834
$ws.Range("M34").Value = -1235.5834

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 76365030
$ws.Range("I11").Value = 739.8333
$ws.Range("K11").Value = 2219.4999
$ws.Range("M11").Value = -2079.4999
$ws.Range("H32").Value = 1040.3334
$ws.Range("I32").Value = 1051
$ws.Range("J32").Value = 1019
$ws.Range("K32").Value = 3153
$ws.Range("L32").Value = 3057
$ws.Range("M32").Value = -2870
$ws.Range("N32").Value = -3623
$ws.Range("H107").Value = 3077.5
$ws.Range("J107").Value = 1893
$ws.Range("L107").Value = 5679
$ws.Range("N107").Value = -9519
$ws.Range("H128").Value = 300000
$ws.Range("I128").Value = 300000
$ws.Range("K128").Value = 900000
$ws.Range("M128").Value = -895020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 55832.55
$ws.Range("I132").Value = 69679.09
$ws.Range("J132").Value = 8138.8887
$ws.Range("K132").Value = 209037.27
$ws.Range("L132").Value = 24416.6661
$ws.Range("M132").Value = -206507.27
$ws.Range("N132").Value = -29476.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6506
$ws.Range("I7").Value = 5650
$ws.Range("J7").Value = 6972.909
$ws.Range("K7").Value = 5650
$ws.Range("L7").Value = 6972.909
$ws.Range("M7").Value = -5538
$ws.Range("N7").Value = -7196.909
$ws.Range("H61").Value = 2487.077
$ws.Range("I61").Value = 1936.1
$ws.Range("K61").Value = 1936.1
$ws.Range("M61").Value = -1734.1
$ws.Range("H68").Value = 3248396.5
$ws.Range("I68").Value = 3789240.2
$ws.Range("J68").Value = 3333
$ws.Range("K68").Value = 3789240.2
$ws.Range("L68").Value = 3333
$ws.Range("M68").Value = -3788491.2
$ws.Range("N68").Value = -4831
$ws.Range("H71").Value = 3248396.5
$ws.Range("I71").Value = 3789240.2
$ws.Range("J71").Value = 3333
$ws.Range("K71").Value = 18946201
$ws.Range("L71").Value = 16665
$ws.Range("M71").Value = -18942457
$ws.Range("N71").Value = -24153
$ws.Range("H113").Value = 2487.077
$ws.Range("I113").Value = 1936.1
$ws.Range("K113").Value = 1936.1
$ws.Range("M113").Value = 233.9000000000001
$ws.Range("H126").Value = 6506
$ws.Range("I126").Value = 5650
$ws.Range("J126").Value = 6972.909
$ws.Range("K126").Value = 16950
$ws.Range("L126").Value = 20918.727
$ws.Range("M126").Value = -14480
$ws.Range("N126").Value = -25858.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 800.125
$ws.Range("I113").Value = 937.55554
$ws.Range("K113").Value = 2812.66662
$ws.Range("M113").Value = -642.66662
$ws.Range("H132").Value = 2320292.5
$ws.Range("I132").Value = 3269101.5
$ws.Range("J132").Value = 16042.714
$ws.Range("K132").Value = 9807304.5
$ws.Range("L132").Value = 48128.142
$ws.Range("M132").Value = -9804774.5
$ws.Range("N132").Value = -53188.142
$ws.Range("H136").Value = 7989.094
$ws.Range("I136").Value = 2265.889
$ws.Range("K136").Value = 6797.667
$ws.Range("M136").Value = -4247.667
